$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '67.487.51'
$ws.Range("E2").Value = '  +1.31%  '

# Row 3
$ws.Range("D3").Value = '3.533.78'
$ws.Range("E3").Value = '  +0.94%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '597.79'
$ws.Range("E5").Value = '  +1.33%  '

# Row 6
Set-TextValue $ws.Range("D6") '173.76'
$ws.Range("E6").Value = '  +2.57%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.592'
$ws.Range("E8").Value = '  +2.06%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.134'
$ws.Range("E9").Value = '  +8.08%  '

# Row 10
$ws.Range("E10").Value = '  +0.81%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.436'
$ws.Range("E11").Value = '  +0.22%  '

# Row 12
$ws.Range("D12").Value = '4.135.64'
$ws.Range("E12").Value = '  +0.74%  '

# Row 13
$ws.Range("E13").Value = '  -0.18%  '

# Row 14
Set-TextValue $ws.Range("D14") '28.79'
$ws.Range("E14").Value = '  +3.01%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.0000182'
$ws.Range("E15").Value = '  +2.88%  '

# Row 16
$ws.Range("D16").Value = '67.345.88'
$ws.Range("E16").Value = '  +1.20%  '

# Row 17
$ws.Range("D17").Value = '3.536.96'
$ws.Range("E17").Value = '  +1.03%  '

# Row 18
Set-TextValue $ws.Range("D18") '6.37'
$ws.Range("E18").Value = '  +1.40%  '

# Row 19
Set-TextValue $ws.Range("D19") '14.21'
$ws.Range("E19").Value = '  +1.68%  '

# Row 20
Set-TextValue $ws.Range("D20") '397.48'
$ws.Range("E20").Value = '  +2.46%  '

# Row 21
Set-TextValue $ws.Range("D21") '8.00'
$ws.Range("E21").Value = '  +0.57%  '

# Row 22
Set-TextValue $ws.Range("D22") '73.53'
$ws.Range("E22").Value = '  +0.79%  '

# Row 23
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D23") '0.541'
$ws.Range("E23").Value = '  +2.81%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D24") '0.999'
$ws.Range("E24").Value = '  -0.19%  '

# Row 25
Set-TextValue $ws.Range("D25") '0.0000124'
$ws.Range("E25").Value = '  +0.71%  '

# Row 26
Set-TextValue $ws.Range("D26") '10.29'
$ws.Range("E26").Value = '  +1.76%  '

# Row 27
Set-TextValue $ws.Range("D27") '0.182'
$ws.Range("E27").Value = '  +0.37%  '

# Row 28
$ws.Range("E28").Value = '  -0.13%  '

# Row 29
Set-TextValue $ws.Range("D29") '6.31'
$ws.Range("E29").Value = '  -0.52%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.47'
$ws.Range("E30").Value = '  +0.21%  '

# Row 31
Set-TextValue $ws.Range("D31") '2.09'
$ws.Range("E31").Value = '  +1.61%  '

# Row 32
Set-TextValue $ws.Range("D32") '24.13'
$ws.Range("E32").Value = '  +3.04%  '

# Row 33
Set-TextValue $ws.Range("D33") '7.40'
$ws.Range("E33").Value = '  +0.00%  '

# Row 34
$ws.Range("E34").Value = '  +4.94%  '

# Row 35
Set-TextValue $ws.Range("D35") '163.90'
$ws.Range("E35").Value = '  +1.67%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.898'
$ws.Range("E36").Value = '  -0.35%  '

# Row 37
Set-TextValue $ws.Range("D37") '1.92'
$ws.Range("E37").Value = '  -0.43%  '

# Row 38
Set-TextValue $ws.Range("D38") '6.97'
$ws.Range("E38").Value = '  +3.96%  '

# Row 39
Set-TextValue $ws.Range("D39") '4.74'
$ws.Range("E39").Value = '  +2.52%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.0748'
$ws.Range("E40").Value = '  +0.25%  '

# Row 41
Set-TextValue $ws.Range("D41") '26.62'
$ws.Range("E41").Value = '  +0.94%  '

# Row 42
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D42") '2.66'
$ws.Range("E42").Value = '  +5.11%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D43") '27.33'
$ws.Range("E43").Value = '  +2.25%  '

# Row 44
$ws.Range("D44").Value = '2.816.25'
$ws.Range("E44").Value = '  +0.81%  '

# Row 45
Set-TextValue $ws.Range("D45") '43.01'
$ws.Range("E45").Value = '  -0.83%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.0311'
$ws.Range("E46").Value = '  -0.89%  '

# Row 47
Set-TextValue $ws.Range("D47") '342.74'
$ws.Range("E47").Value = '  -2.66%  '

# Row 48
$ws.Range("E48").Value = '  +0.93%  '

# Row 49
Set-TextValue $ws.Range("D49") '33.97'
$ws.Range("E49").Value = '  +2.54%  '

# Row 50
Set-TextValue $ws.Range("D50") '6.53'
$ws.Range("E50").Value = '  +0.89%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.855'
$ws.Range("E51").Value = '  +0.72%  '
